$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.253.26"
$ws.Range("E2").Value = "  -1.86%  "
$ws.Range("D3").Value = "'2.751.24"
$ws.Range("E3").Value = "  -2.90%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'355.16"
$ws.Range("E5").Value = "  -1.96%  "
$ws.Range("D6").Value = "'107.81"
$ws.Range("E6").Value = "  -3.80%  "
$ws.Range("E7").Value = "  -3.13%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "'0.578"
$ws.Range("E9").Value = "  -4.26%  "
$ws.Range("D10").Value = "'39.18"
$ws.Range("E10").Value = "  -4.32%  "
$ws.Range("E11").Value = "  +3.17%  "
$ws.Range("D12").Value = "'0.0834"
$ws.Range("E12").Value = "  -3.68%  "
$ws.Range("D13").Value = "'19.53"
$ws.Range("E13").Value = "  -2.38%  "
$ws.Range("D14").Value = "'7.49"
$ws.Range("E14").Value = "  -4.07%  "
$ws.Range("D15").Value = "'3.192.56"
$ws.Range("E15").Value = "  -2.82%  "
$ws.Range("D16").Value = "'2.744.79"
$ws.Range("E16").Value = "  -2.98%  "
$ws.Range("D17").Value = "'0.920"
$ws.Range("E17").Value = "  -1.32%  "
$ws.Range("D18").Value = "'51.166.88"
$ws.Range("E18").Value = "  -1.79%  "
$ws.Range("D19").Value = "'7.51"
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("E20").Value = "  -4.59%  "
$ws.Range("D21").Value = "'12.98"
$ws.Range("E21").Value = "  -3.25%  "
$ws.Range("D22").Value = "'0.0₃0957"
$ws.Range("E22").Value = "  -4.08%  "
$ws.Range("D23").Value = "'69.43"
$ws.Range("E23").Value = "  -1.46%  "
$ws.Range("D24").Value = "'264.39"
$ws.Range("E24").Value = "  -3.00%  "
$ws.Range("D25").Value = "'2.72"
$ws.Range("E25").Value = "  -3.53%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "'25.94"
$ws.Range("E27").Value = "  -3.56%  "
$ws.Range("E28").Value = "  +14.47%  "
$ws.Range("D29").Value = "'10.07"
$ws.Range("E29").Value = "  -2.46%  "
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("D31").Value = "'34.89"
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'6.06"
$ws.Range("E32").Value = "  +3.02%  "
$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").Value = "'51.52"
$ws.Range("E33").Value = "  -1.53%  "
$ws.Range("D34").Value = "'0.0439"
$ws.Range("E34").Value = "  -8.19%  "
$ws.Range("D35").Value = "'0.0830"
$ws.Range("E35").Value = "  -2.62%  "
$ws.Range("D36").Value = "'5.13"
$ws.Range("E36").Value = "  -8.01%  "
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("D38").Value = "'18.50"
$ws.Range("E38").Value = "  +0.74%  "
$ws.Range("D39").Value = "'3.12"
$ws.Range("E39").Value = "  -4.72%  "
$ws.Range("D40").Value = "'1.93"
$ws.Range("E40").Value = "  -5.21%  "
$ws.Range("D41").Value = "'0.113"
$ws.Range("E41").Value = "  -3.14%  "
$ws.Range("D42").Value = "'2.48"
$ws.Range("E42").Value = "  -2.55%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "'2.21"
$ws.Range("E43").Value = "  -2.61%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").Value = "'118.86"
$ws.Range("E44").Value = "  -4.67%  "
$ws.Range("E45").Value = "  -6.61%  "
$ws.Range("D46").Value = "'2.074.53"
$ws.Range("E46").Value = "  -0.51%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'3.21"
$ws.Range("E47").Value = "  -3.55%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "'2.30"
$ws.Range("E48").Value = "  -1.22%  "
$ws.Range("B49").Value = "SEI"
$ws.Range("C49").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
$ws.Range("D49").Value = "'0.919"
$ws.Range("E49").Value = "  -4.84%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").Value = "'5.48"
$ws.Range("E50").Value = "  -6.79%  "
$ws.Range("E51").Value = "  +3.90%  "
